$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The edit: in the "Количество игроков : 2" paragraph, drop the
# gramStart/gramEnd-wrapped " :" / " 2" runs and replace them with a
# single run containing ": 2"; the _GoBack bookmark (previously sitting
# right before the "В целом, правильная игра..." paragraph) is moved
# so that it now sits right after the "Количество игроков" run, before
# the new ": 2" run.
# ------------------------------------------------------------------

# Locate "Количество игроков" in the document.
$findRange = $d.Content
$findRange.Find.Execute("Количество игроков", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$labelStart = $findRange.Start
$labelLen = "Количество игроков".Length

# Range covering "Количество игроков : 2" (label + " : 2" that follows it).
$wholeRange = $d.Range($labelStart, $labelStart + $labelLen + 4)

# Re-write the whole chunk as a single clean run, dropping the
# gramStart/gramEnd proofErr marks and the extra space before the colon.
$wholeRange.Text = "Количество игроков: 2"

# Insert (move) the _GoBack bookmark right after "Количество игроков",
# i.e. right before ": 2". Word only keeps a single _GoBack bookmark, so
# adding it here automatically removes it from its old location.
$bmPoint = $d.Range($labelStart + $labelLen, $labelStart + $labelLen)
$d.Bookmarks.Add("_GoBack", $bmPoint)
